# Commit: "Added target species CALLSAP in excel file"
# The target_species sheet lists one species per row (species_name, target)
# with the used range currently ending at row 31 (GALEECH). Append the new
# target species "CALLSAP" as the next row, flagged as a target (1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "CALLSAP"
$ws.Cells.Item($newRow, 2).Value = 1
